# Auto-generated edit script: Add data for 2024-03-13
# Updates violent crime counts (mostly column K = year 2024 partial-year total,
# plus a few corrections to H/I for 2021/2022) across the Citywide Totals,
# By Neighborhood, and individual neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 1379
$ws.Range('K3').Value = 1312
$ws.Range('H4').Value = 1720
$ws.Range('K4').Value = 287
$ws.Range('I5').Value = 723
$ws.Range('K5').Value = 85
$ws.Range('K6').Value = 1673
$ws.Range('H7').Value = 26031
$ws.Range('I7').Value = 26241
$ws.Range('K7').Value = 4736

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 66

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 84
$ws.Range('K3').Value = 81
$ws.Range('K6').Value = 90
$ws.Range('K7').Value = 276

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 44
$ws.Range('K7').Value = 96

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 55
$ws.Range('K3').Value = 77
$ws.Range('K4').Value = 11
$ws.Range('K7').Value = 195

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K3').Value = 23
$ws.Range('K5').Value = 2
$ws.Range('K7').Value = 81

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K5').Value = 9
$ws.Range('K6').Value = 59
$ws.Range('K7').Value = 160

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 32
$ws.Range('K7').Value = 120

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 24
$ws.Range('K6').Value = 26
$ws.Range('K7').Value = 86

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K6').Value = 39
$ws.Range('K7').Value = 139
$ws.Range('K8').Value = 276
$ws.Range('K18').Value = 37
$ws.Range('K19').Value = 125
$ws.Range('K20').Value = 112
$ws.Range('K23').Value = 46
$ws.Range('K24').Value = 17
$ws.Range('K29').Value = 218
$ws.Range('K31').Value = 52
$ws.Range('K33').Value = 195
$ws.Range('K36').Value = 53
$ws.Range('K37').Value = 160
$ws.Range('K41').Value = 50
$ws.Range('K42').Value = 162
$ws.Range('K43').Value = 46
$ws.Range('K47').Value = 35
$ws.Range('K48').Value = 52
$ws.Range('K49').Value = 32
$ws.Range('K52').Value = 126
$ws.Range('K53').Value = 66
$ws.Range('K54').Value = 83
$ws.Range('K55').Value = 50
$ws.Range('K60').Value = 35
$ws.Range('H63').Value = 273
$ws.Range('I63').Value = 197
$ws.Range('K63').Value = 17
$ws.Range('K64').Value = 31
$ws.Range('K65').Value = 120
$ws.Range('K67').Value = 189
$ws.Range('K76').Value = 65
$ws.Range('K77').Value = 35
$ws.Range('K78').Value = 65
$ws.Range('K79').Value = 128
$ws.Range('K80').Value = 19
$ws.Range('K83').Value = 96
$ws.Range('K84').Value = 34
$ws.Range('K85').Value = 245
$ws.Range('K86').Value = 34
$ws.Range('K88').Value = 59
$ws.Range('K89').Value = 64
$ws.Range('K90').Value = 42
$ws.Range('K92').Value = 20
$ws.Range('K93').Value = 21
$ws.Range('K94').Value = 59
$ws.Range('K95').Value = 81
$ws.Range('K97').Value = 41
$ws.Range('K99').Value = 86
$ws.Range('H101').Value = 26031
$ws.Range('I101').Value = 26241
$ws.Range('K101').Value = 4736

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 52

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 55
$ws.Range('K3').Value = 56
$ws.Range('K7').Value = 189

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K6').Value = 10
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 32

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 27
$ws.Range('K6').Value = 33
$ws.Range('K7').Value = 83

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 70
$ws.Range('K7').Value = 218

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 52

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 41
$ws.Range('K3').Value = 38
$ws.Range('K7').Value = 125

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 65

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 13
$ws.Range('K3').Value = 13
$ws.Range('K6').Value = 12
$ws.Range('K7').Value = 39

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K3').Value = 9
$ws.Range('K7').Value = 50

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 47
$ws.Range('K5').Value = 2
$ws.Range('K6').Value = 64
$ws.Range('K7').Value = 162

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K5').Value = 2
$ws.Range('K7').Value = 65

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 17
$ws.Range('K7').Value = 50

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 17

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 18
$ws.Range('K3').Value = 13
$ws.Range('K7').Value = 46

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 47
$ws.Range('K3').Value = 43
$ws.Range('K6').Value = 27
$ws.Range('K7').Value = 128

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 9
$ws.Range('K7').Value = 31

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 33
$ws.Range('K6').Value = 43
$ws.Range('K7').Value = 112

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K3').Value = 8
$ws.Range('K7').Value = 37

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 20
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 53

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K3').Value = 5
$ws.Range('K7').Value = 21

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 53
$ws.Range('K3').Value = 42
$ws.Range('K6').Value = 35
$ws.Range('K7').Value = 139

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K2').Value = 19
$ws.Range('K6').Value = 26
$ws.Range('K7').Value = 59

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K3').Value = 12
$ws.Range('K7').Value = 35

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 27
$ws.Range('K7').Value = 41

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('K2').Value = 5
$ws.Range('K7').Value = 20

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K6').Value = 33
$ws.Range('K7').Value = 59

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 9
$ws.Range('K7').Value = 64

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K3').Value = 6
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K6').Value = 11
$ws.Range('K7').Value = 42

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K2').Value = 5
$ws.Range('K3').Value = 15
$ws.Range('K7').Value = 35

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K2').Value = 8
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 46

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 89
$ws.Range('K3').Value = 78
$ws.Range('K6').Value = 61
$ws.Range('K7').Value = 245

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K3').Value = 13
$ws.Range('K7').Value = 35

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K2').Value = 4
$ws.Range('K7').Value = 19

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 28
$ws.Range('K7').Value = 126
